$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 9.031965666666666
$ws.Range("H2").Value = 27.095897
$ws.Range("I2").Value = 0.4424406034784756
$ws.Range("J2").Value = 0.4424406034784755
$ws.Range("M2").Value = 4.066173333333333
$ws.Range("N2").Value = 12.19852
$ws.Range("O2").Value = 0.8070107842953054
$ws.Range("P2").Value = 0.8070107842953055
$ws.Range("Q2").Value = 36.72553794138221
$ws.Range("R2").Value = 330.5298414724399
$ws.Range("S2").Value = 0.3570543384172528
$ws.Range("T2").Value = 0.3570543384172528

# Row 3
$ws.Range("G3").Value = 9.031965666666666
$ws.Range("H3").Value = 27.095897
$ws.Range("I3").Value = 0.4424406034784756
$ws.Range("J3").Value = 0.4424406034784755
$ws.Range("O3").Value = 0.03600419273120554
$ws.Range("P3").Value = 0.03600419273120554
$ws.Range("Q3").Value = 1.638482870279555
$ws.Range("R3").Value = 14.746345832516
$ws.Range("S3").Value = 0.01592971675974993
$ws.Range("T3").Value = 0.01592971675974992

# Row 4
$ws.Range("G4").Value = 9.031965666666666
$ws.Range("H4").Value = 27.095897
$ws.Range("I4").Value = 0.4424406034784756
$ws.Range("J4").Value = 0.4424406034784755
$ws.Range("M4").Value = 0.7909786666666667
$ws.Range("N4").Value = 2.372936
$ws.Range("O4").Value = 0.156985022973489
$ws.Range("P4").Value = 0.156985022973489
$ws.Range("Q4").Value = 7.144092160399111
$ws.Range("R4").Value = 64.296829443592
$ws.Range("S4").Value = 0.06945654830147283
$ws.Range("T4").Value = 0.06945654830147283

# Row 5
$ws.Range("I5").Value = 0.4469933372071527
$ws.Range("J5").Value = 0.4469933372071526
$ws.Range("M5").Value = 4.066173333333333
$ws.Range("N5").Value = 12.19852
$ws.Range("O5").Value = 0.8070107842953054
$ws.Range("P5").Value = 0.8070107842953055
$ws.Range("Q5").Value = 37.10344538019999
$ws.Range("R5").Value = 333.9310084217999
$ws.Range("S5").Value = 0.3607284436343202
$ws.Range("T5").Value = 0.3607284436343202

# Row 6
$ws.Range("I6").Value = 0.4469933372071527
$ws.Range("J6").Value = 0.4469933372071526
$ws.Range("O6").Value = 0.03600419273120554
$ws.Range("P6").Value = 0.03600419273120554
$ws.Range("S6").Value = 0.01609363426237108
$ws.Range("T6").Value = 0.01609363426237107

# Row 7
$ws.Range("I7").Value = 0.4469933372071527
$ws.Range("J7").Value = 0.4469933372071526
$ws.Range("M7").Value = 0.7909786666666667
$ws.Range("N7").Value = 2.372936
$ws.Range("O7").Value = 0.156985022973489
$ws.Range("P7").Value = 0.156985022973489
$ws.Range("Q7").Value = 7.217605190359999
$ws.Range("R7").Value = 64.95844671324
$ws.Range("S7").Value = 0.07017125931046138
$ws.Range("T7").Value = 0.07017125931046138

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.2635683333333333
$ws.Range("H8").Value = 0.790705
$ws.Range("I8").Value = 0.01291117977653399
$ws.Range("J8").Value = 0.01291117977653399
$ws.Range("M8").Value = 4.066173333333333
$ws.Range("N8").Value = 12.19852
$ws.Range("O8").Value = 0.8070107842953054
$ws.Range("P8").Value = 0.8070107842953055
$ws.Range("Q8").Value = 1.071714528511111
$ws.Range("R8").Value = 9.645430756599998
$ws.Range("S8").Value = 0.01041946131763838
$ws.Range("T8").Value = 0.01041946131763838

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.2635683333333333
$ws.Range("H9").Value = 0.790705
$ws.Range("I9").Value = 0.01291117977653399
$ws.Range("J9").Value = 0.01291117977653399
$ws.Range("O9").Value = 0.03600419273120554
$ws.Range("P9").Value = 0.03600419273120554
$ws.Range("Q9").Value = 0.04781375563777777
$ws.Range("R9").Value = 0.4303238007399999
$ws.Range("S9").Value = 0.0004648566050615733
$ws.Range("T9").Value = 0.0004648566050615731

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2635683333333333
$ws.Range("H10").Value = 0.790705
$ws.Range("I10").Value = 0.01291117977653399
$ws.Range("J10").Value = 0.01291117977653399
$ws.Range("M10").Value = 0.7909786666666667
$ws.Range("N10").Value = 2.372936
$ws.Range("O10").Value = 0.156985022973489
$ws.Range("P10").Value = 0.156985022973489
$ws.Range("Q10").Value = 0.2084769288755556
$ws.Range("R10").Value = 1.87629235988
$ws.Range("S10").Value = 0.002026861853834036
$ws.Range("T10").Value = 0.002026861853834035

# Row 11
$ws.Range("G11").Value = 1.809602666666667
$ws.Range("H11").Value = 5.428808
$ws.Range("I11").Value = 0.088645343156153
$ws.Range("J11").Value = 0.08864534315615297
$ws.Range("M11").Value = 4.066173333333333
$ws.Range("N11").Value = 12.19852
$ws.Range("O11").Value = 0.8070107842953054
$ws.Range("P11").Value = 0.8070107842953055
$ws.Range("Q11").Value = 7.358158107128888
$ws.Range("R11").Value = 66.22342296416
$ws.Range("S11").Value = 0.07153774790457351
$ws.Range("T11").Value = 0.0715377479045735

# Row 12
$ws.Range("G12").Value = 1.809602666666667
$ws.Range("H12").Value = 5.428808
$ws.Range("I12").Value = 0.088645343156153
$ws.Range("J12").Value = 0.08864534315615297
$ws.Range("O12").Value = 0.03600419273120554
$ws.Range("P12").Value = 0.03600419273120554
$ws.Range("Q12").Value = 0.3282788133582222
$ws.Range("R12").Value = 2.954509320224
$ws.Range("S12").Value = 0.003191604019717985
$ws.Range("T12").Value = 0.003191604019717984

# Row 13
$ws.Range("G13").Value = 1.809602666666667
$ws.Range("H13").Value = 5.428808
$ws.Range("I13").Value = 0.088645343156153
$ws.Range("J13").Value = 0.08864534315615297
$ws.Range("M13").Value = 0.7909786666666667
$ws.Range("N13").Value = 2.372936
$ws.Range("O13").Value = 0.156985022973489
$ws.Range("P13").Value = 0.156985022973489
$ws.Range("Q13").Value = 1.431357104476445
$ws.Range("R13").Value = 12.882213940288
$ws.Range("S13").Value = 0.01391599123186149
$ws.Range("T13").Value = 0.01391599123186149

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1839203333333334
$ws.Range("H14").Value = 0.5517610000000001
$ws.Range("I14").Value = 0.009009536381684918
$ws.Range("J14").Value = 0.009009536381684917
$ws.Range("M14").Value = 4.066173333333333
$ws.Range("N14").Value = 12.19852
$ws.Range("O14").Value = 0.8070107842953054
$ws.Range("P14").Value = 0.8070107842953055
$ws.Range("Q14").Value = 0.7478519548577777
$ws.Range("R14").Value = 6.73066759372
$ws.Range("S14").Value = 0.007270793021520634
$ws.Range("T14").Value = 0.007270793021520634

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1839203333333334
$ws.Range("H15").Value = 0.5517610000000001
$ws.Range("I15").Value = 0.009009536381684918
$ws.Range("J15").Value = 0.009009536381684917
$ws.Range("O15").Value = 0.03600419273120554
$ws.Range("P15").Value = 0.03600419273120554
$ws.Range("Q15").Value = 0.03336486505644445
$ws.Range("R15").Value = 0.300283785508
$ws.Range("S15").Value = 0.000324381084304992
$ws.Range("T15").Value = 0.000324381084304992

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1839203333333334
$ws.Range("H16").Value = 0.5517610000000001
$ws.Range("I16").Value = 0.009009536381684918
$ws.Range("J16").Value = 0.009009536381684917
$ws.Range("M16").Value = 0.7909786666666667
$ws.Range("N16").Value = 2.372936
$ws.Range("O16").Value = 0.156985022973489
$ws.Range("P16").Value = 0.156985022973489
$ws.Range("Q16").Value = 0.1454770600328889
$ws.Range("R16").Value = 1.309293540296
$ws.Range("S16").Value = 0.001414362275859292
$ws.Range("T16").Value = 0.001414362275859292
